$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

# 1. MVPA quartile boundary labels (ENMO column)
Replace-Text "[0.917,320]" "[1.83,319]"
Replace-Text "(320,464]" "(319,464]"

# 2. Row (319/320,464] HR / p-value (ENMO column)
Replace-Text "0.69 (0.54 to 0.88)" "0.67 (0.52 to 0.87)"

# 3. Row (464,642] HR / p-value (ENMO column)
Replace-Text "0.75 (0.58 to 0.96)" "0.70 (0.54 to 0.92)"
Replace-Text "0.023" "0.010"

# 4. Row (642,2.39e+03] HR / p-value (ENMO column)
Replace-Text "0.59 (0.44 to 0.78)" "0.61 (0.45 to 0.82)"
Replace-Text "<0.001" "0.001"

# 5. Row (115,234] HR / p-value (Random Forest and HMM column)
Replace-Text "0.82 (0.65 to 1.03)" "0.81 (0.63 to 1.04)"
Replace-Text "0.094" "0.095"

# 6. Row (234,404] HR / p-value (Random Forest and HMM column)
Replace-Text "0.73 (0.56 to 0.94)" "0.72 (0.55 to 0.95)"
Replace-Text "0.017" "0.022"

# 7. Row (404,2.47e+03] HR / p-value (Random Forest and HMM column)
Replace-Text "0.76 (0.57 to 1.02)" "0.83 (0.61 to 1.12)"
Replace-Text "0.069" "0.22"
